# Auto-generated edit script: updates crypto price/volume table
# to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.085.78'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '2.016.97'
$ws.Range('E3').Value = '  -1.47%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.993'
$ws.Range('E4').Value = '  -0.95%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.34'
$ws.Range('E5').Value = '  -1.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.608'
$ws.Range('E6').Value = '  -1.76%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '54.80'
$ws.Range('E8').Value = '  -3.58%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.380'
$ws.Range('E9').Value = '  -0.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0784'
$ws.Range('E10').Value = '  +2.17%  '
$ws.Range('E11').Value = '  -2.97%  '
$ws.Range('D12').Value = '2.315.49'
$ws.Range('E12').Value = '  -1.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.30'
$ws.Range('E13').Value = '  -2.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.40'
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.740'
$ws.Range('E15').Value = '  -1.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.13'
$ws.Range('E16').Value = '  -1.93%  '
$ws.Range('D17').Value = '2.017.42'
$ws.Range('E17').Value = '  -1.54%  '
$ws.Range('D18').Value = '36.981.82'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.15'
$ws.Range('E19').Value = '  +3.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.84'
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('D21').Value = '0.0₃0821'
$ws.Range('E21').Value = '  +0.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '226.24'
$ws.Range('E22').Value = '  +0.39%  '
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.41'
$ws.Range('E24').Value = '  +3.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.19'
$ws.Range('E25').Value = '  -5.46%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.31'
$ws.Range('E26').Value = '  -1.73%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.21'
$ws.Range('E27').Value = '  -2.61%  '
$ws.Range('E28').Value = '  -1.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.74'
$ws.Range('E29').Value = '  -1.64%  '
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('E31').Value = '  -3.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.51'
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0616'
$ws.Range('E33').Value = '  -0.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.42'
$ws.Range('E34').Value = '  -3.19%  '
$ws.Range('E35').Value = '  -3.98%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.84'
$ws.Range('E36').Value = '  +1.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.15'
$ws.Range('E38').Value = '  -3.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.35'
$ws.Range('E40').Value = '  -3.52%  '
$ws.Range('D41').Value = '1.478.45'
$ws.Range('E41').Value = '  -0.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '16.75'
$ws.Range('E42').Value = '  +1.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '95.34'
$ws.Range('E43').Value = '  -2.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0921'
$ws.Range('E44').Value = '  -2.73%  '
$ws.Range('E45').Value = '  -4.12%  '
$ws.Range('E46').Value = '  -3.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.24'
$ws.Range('E47').Value = '  +1.32%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.01'
$ws.Range('E48').Value = '  -1.59%  '
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('D50').Value = '2.201.88'
$ws.Range('E50').Value = '  -1.42%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.62'
$ws.Range('E51').Value = '  -10.07%  '
